$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.190686464309692
$ws.Range("B1").Value = 2.473435401916504
$ws.Range("D1").Value = 2.278876066207886
$ws.Range("E1").Value = 1.179906010627747
